# feat: added student profile
# Adds new student-profile columns (D:L) to the header row of Sheet1 and
# selects the full header row, mirroring the authored OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value2 = "contactNo"
$ws.Range("E1").Value2 = "studentTenthMarks"
$ws.Range("F1").Value2 = "studentTwelthMarks"
$ws.Range("G1").Value2 = "studentUGMarks"
$ws.Range("H1").Value2 = "studentPGMarks"
$ws.Range("I1").Value2 = "studentDescription"
$ws.Range("J1").Value2 = "studentId"
$ws.Range("K1").Value2 = "dept"
$ws.Range("L1").Value2 = "gender"

$ws.Rows("1:1").Select()
